# "Added last minute updates"
#
# The document's first paragraph is a placeholder/bookmark line:
#   **ID__AFFARS_pgi_5301_topic_31__ID**<space>
# (two runs: the marker text, and a trailing literal space run).
#
# This edit:
#   1. Renames the placeholder marker to **ID__AFFARS_AFMC_PGI_5301_90__ID**
#      and drops the now-superfluous trailing-space run, leaving a single run.
#   2. Gives the paragraph a (currently invisible/no-line) paragraph border
#      with 5-twip spacing on all four sides, matching the border already
#      used further down the document.
#   3. Widens the paragraph's left indent from 120 twips (6pt) to
#      225 twips (11.25pt), again matching the sibling paragraph below it.

$d = $word.ActiveDocument

$p = $d.Paragraphs.Item(1)
$pRange = $p.Range

# Locate the marker text inside paragraph 1.
$markerRange = $pRange.Duplicate
$found = $markerRange.Find.Execute("**ID__AFFARS_pgi_5301_topic_31__ID**", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)

# Remove anything left in the paragraph after the marker (the trailing-space
# run) but keep the paragraph mark itself.
$trailing = $d.Range($markerRange.End, $pRange.End - 1)
if ($trailing.Start -lt $trailing.End) {
    $trailing.Delete()
}

# Rename the marker itself (this preserves the run's existing rPr).
$markerRange.Text = "**ID__AFFARS_AFMC_PGI_5301_90__ID**"

# Apply the new paragraph formatting (border + indent) to paragraph 1.
$pf = $d.Paragraphs.Item(1).Range.ParagraphFormat
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 11.25
